$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the existing data rows (2-8) before any structural changes.
$oldA5 = $ws.Range("A5").Value2
$oldB5 = $ws.Range("B5").Value2
$oldC5 = $ws.Range("C5").Value2

$oldA6 = $ws.Range("A6").Value2
$oldB6 = $ws.Range("B6").Value2
$oldC6 = $ws.Range("C6").Value2

$oldA7 = $ws.Range("A7").Value2
$oldB7 = $ws.Range("B7").Value2
$oldC7 = $ws.Range("C7").Value2

$oldA8 = $ws.Range("A8").Value2
$oldB8 = $ws.Range("B8").Value2
$oldC8 = $ws.Range("C8").Value2

# Insert a new row at position 5 (pushes existing rows 5-8 down to 6-9)
$ws.Rows.Item(5).Insert()
# Insert a new row at position 10 (after the now-existing row 9)
$ws.Rows.Item(10).Insert()

# Fill in the new rows' "vernacularName" (column B) values first
$ws.Range("B5").Value2 = "五色鳥"
$ws.Range("B10").Value2 = "臺灣海棗"

# Then the new rows' "storageLocation" (column C) values
$ws.Range("C5").Value2 = "A4"
$ws.Range("C10").Value2 = "P3"

# Then the new rows' "occurrenceID" (column A) values
$ws.Range("A5").Value2 = "bb514b5d-d30a-42e8-bf56-fe58063a5892"
$ws.Range("A10").Value2 = "04b93e67-7389-435b-9113-2936c4e1f3d3"

# Re-write the pre-existing taxa into their new row positions:
#  old row8 -> new row6, old row5 -> new row7, old row6 -> new row8, old row7 -> new row9
$ws.Range("A6").Value2 = $oldA8
$ws.Range("B6").Value2 = $oldB8
$ws.Range("C6").Value2 = $oldC8

$ws.Range("A7").Value2 = $oldA5
$ws.Range("B7").Value2 = $oldB5
$ws.Range("C7").Value2 = $oldC5

$ws.Range("A8").Value2 = $oldA6
$ws.Range("B8").Value2 = $oldB6
$ws.Range("C8").Value2 = $oldC6

$ws.Range("A9").Value2 = $oldA7
$ws.Range("B9").Value2 = $oldB7
$ws.Range("C9").Value2 = $oldC7

# Fix up column-A formatting: rows 5, 6, 9, 10 carry no explicit cell style,
# while rows 7 and 8 keep the style used throughout the rest of the table.
$ws.Range("A5").ClearFormats()
$ws.Range("A6").ClearFormats()
$ws.Range("A9").ClearFormats()
$ws.Range("A10").ClearFormats()

$ws.Range("A2").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

Write-Output "done"
